$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.347.90"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.933.22"
$ws.Range("E3").Value = "  -2.13%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'241.33"
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("D6").Value = "'0.606"
$ws.Range("E6").Value = "  -2.75%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'56.36"
$ws.Range("E8").Value = "  -3.62%  "
$ws.Range("E9").Value = "  -3.74%  "
$ws.Range("D10").Value = "'0.0832"
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("E11").Value = "  -1.66%  "
$ws.Range("D12").Value = "2.214.04"
$ws.Range("E12").Value = "  -2.30%  "
$ws.Range("D13").Value = "'21.15"
$ws.Range("E13").Value = "  -7.64%  "
$ws.Range("D14").Value = "'0.798"
$ws.Range("E14").Value = "  -6.54%  "
$ws.Range("D15").Value = "'13.29"
$ws.Range("E15").Value = "  -4.17%  "
$ws.Range("D16").Value = "'5.11"
$ws.Range("E16").Value = "  -5.54%  "
$ws.Range("D17").Value = "1.938.99"
$ws.Range("E17").Value = "  -1.71%  "
$ws.Range("D18").Value = "36.294.13"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").Value = "0.0₃0858"
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("D20").Value = "'68.56"
$ws.Range("E20").Value = "  -2.37%  "
$ws.Range("D21").Value = "'225.92"
$ws.Range("E21").Value = "  -3.27%  "
$ws.Range("E22").Value = "  -6.27%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'2.33"
$ws.Range("E24").Value = "  -7.18%  "
$ws.Range("D25").Value = "'2.26"
$ws.Range("E25").Value = "  -1.00%  "
$ws.Range("D26").Value = "'9.04"
$ws.Range("E26").Value = "  -7.89%  "
$ws.Range("D27").Value = "'160.55"
$ws.Range("E27").Value = "  -1.52%  "
$ws.Range("E28").Value = "  -2.53%  "
$ws.Range("D29").Value = "'19.05"
$ws.Range("E29").Value = "  -3.35%  "
$ws.Range("D30").Value = "'0.117"
$ws.Range("E30").Value = "  -2.18%  "
$ws.Range("E31").Value = "  -6.03%  "
$ws.Range("D32").Value = "'4.50"
$ws.Range("E32").Value = "  -7.01%  "
$ws.Range("D33").Value = "'0.0619"
$ws.Range("E33").Value = "  -8.98%  "
$ws.Range("D34").Value = "'4.13"
$ws.Range("E34").Value = "  -5.62%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  -1.62%  "
$ws.Range("D37").Value = "'5.92"
$ws.Range("E37").Value = "  -3.96%  "
$ws.Range("E38").Value = "  -3.87%  "
$ws.Range("D39").Value = "'2.96"
$ws.Range("E39").Value = "  +1.19%  "
$ws.Range("D40").Value = "'0.0964"
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("D41").Value = "'2.86"
$ws.Range("E41").Value = "  -1.01%  "
$ws.Range("E42").Value = "  -2.53%  "
$ws.Range("E43").Value = "  -6.77%  "
$ws.Range("D44").Value = "'15.56"
$ws.Range("E44").Value = "  -2.57%  "
$ws.Range("D45").Value = "1.321.80"
$ws.Range("E45").Value = "  -2.64%  "
$ws.Range("E46").Value = "  -6.82%  "
$ws.Range("D47").Value = "'84.57"
$ws.Range("E47").Value = "  -7.53%  "
$ws.Range("E48").Value = "  -6.18%  "
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").Value = "2.105.54"
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("D51").Value = "'43.04"
$ws.Range("E51").Value = "  -3.67%  "
